$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - these values may look numeric, so we force them
# to remain text (matching the original inlineStr cell type) by using the
# classic leading-apostrophe text-entry marker, then restoring the cell's
# original Style afterward so no stray number-format is left applied to the cell.
$priceUpdates = @{
    'D2' = '29.501.02'
    'D3' = '1.857.76'
    'D4' = '0.9995'
    'D5' = '245.40'
    'D6' = '0.6969'
    'D8' = '0.3078'
    'D9' = '0.07703'
    'D10' = '23.68'
    'D11' = '0.07793'
    'D12' = '5.161'
    'D13' = '1.853.53'
    'D14' = '0.6950'
    'D15' = '91.25'
    'D16' = '6.339'
    'D17' = '29.483.28'
    'D18' = '0.000008321'
    'D19' = '2.101.15'
    'D20' = '238.63'
    'D21' = '12.75'
    'D22' = '1.000'
    'D23' = '7.629'
    'D26' = '160.17'
    'D27' = '8.902'
    'D28' = '18.29'
    'D29' = '1.535'
    'D30' = '4.251'
    'D31' = '4.152'
    'D32' = '1.205'
    'D33' = '0.05110'
    'D34' = '0.7761'
    'D35' = '1.883'
    'D37' = '2.686'
    'D38' = '1.318.53'
    'D39' = '0.01877'
    'D40' = '2.726'
    'D41' = '0.9528'
    'D42' = '106.02'
    'D43' = '5.777'
    'D45' = '9.810'
    'D47' = '2.000.90'
    'D49' = '1.789'
    'D50' = '63.18'
    'D51' = '6.967'
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.Value = "'" + $priceUpdates[$addr]
    $cell.Style = $origStyle
}

# Volume(1h) (column E) updates - these are always non-numeric (contain "%"
# and surrounding spaces) so a plain Value assignment keeps them as text.
$volumeUpdates = @{
    'E2' = '  +2.10%  '
    'E3' = '  +1.37%  '
    'E4' = '  -0.04%  '
    'E5' = '  +0.21%  '
    'E6' = '  +1.16%  '
    'E7' = '  -0.02%  '
    'E8' = '  +0.93%  '
    'E9' = '  +0.18%  '
    'E10' = '  +1.26%  '
    'E11' = '  -0.17%  '
    'E12' = '  +1.52%  '
    'E13' = '  +1.09%  '
    'E14' = '  +1.98%  '
    'E15' = '  +0.94%  '
    'E17' = '  +2.07%  '
    'E18' = '  +0.22%  '
    'E19' = '  +1.31%  '
    'E20' = '  -1.57%  '
    'E21' = '  +0.27%  '
    'E22' = '  -0.02%  '
    'E23' = '  +2.27%  '
    'E24' = '  -0.01%  '
    'E25' = '  +1.29%  '
    'E26' = '  -0.63%  '
    'E27' = '  +1.06%  '
    'E28' = '  +0.58%  '
    'E29' = '  -0.58%  '
    'E30' = '  +0.91%  '
    'E31' = '  +0.00%  '
    'E32' = '  +1.87%  '
    'E33' = '  -0.06%  '
    'E34' = '  +1.36%  '
    'E35' = '  +2.30%  '
    'E36' = '  +0.78%  '
    'E37' = '  -0.36%  '
    'E38' = '  +8.30%  '
    'E39' = '  +1.76%  '
    'E40' = '  +1.08%  '
    'E41' = '  +1.37%  '
    'E42' = '  -2.51%  '
    'E43' = '  +0.85%  '
    'E44' = '  +0.12%  '
    'E45' = '  +3.06%  '
    'E46' = '  +1.88%  '
    'E47' = '  +1.32%  '
    'E48' = '  +1.35%  '
    'E49' = '  +2.37%  '
    'E50' = '  -1.58%  '
    'E51' = '  +1.12%  '
}

foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
